# Insert one new weekly record as row 97 on the single sheet, pushing the
# existing rows 97..163 down to 98..164 (dimension A1:T163 -> A1:T164).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(97).Insert()

$ws.Range("A97").Value = 9
$ws.Range("B97").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C97").Value = "Metropolitana"
$ws.Range("D97").Value = 45062
$ws.Range("E97").Value = 13
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100101004
$ws.Range("J97").Value = "Frambuesa"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 470
$ws.Range("N97").Value = 7500
$ws.Range("O97").Value = 8000
$ws.Range("P97").Value = 7734
$ws.Range("Q97").Value = "$/bandeja 2 kilos"
$ws.Range("R97").Value = "Provincia de Linares"
$ws.Range("S97").Value = 3867
$ws.Range("T97").Value = 2
